# ---------------------------------------------------------------------------
# Apply the "Read Data from the Excel Sheet" update:
#   - Admin sheet (sheet1) and PIM sheet (sheet2) get a row of test-case
#     metadata (name / userName / passWord / searchUser / name-repeated),
#     each cell boxed with a thin border; the first/last cell of each row
#     additionally gets a yellow highlight fill.
#   - A brand-new "TIME" worksheet is appended with four repetitions of the
#     same boxed row pattern (some rows carry two extra logged columns).
#   - Column widths are autofit, selections + the active tab are updated.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$YELLOW = 65535   # RGB(255,255,0)
$xlEdgeRight = 10

# Paint a thin box border (all four edges) around $rng, then optionally
# punch out the right edge and/or apply the yellow highlight fill. Doing the
# border first and the fill second (in that order) lets matching cells reuse
# the same style record instead of minting a new one for every call.
function Format-Cell {
    param($rng, [bool]$dropRight = $false, [bool]$fill = $false)
    $rng.Borders.LineStyle = 1
    if ($dropRight) {
        $rng.Borders.Item($xlEdgeRight).LineStyle = 0
    }
    if ($fill) {
        $rng.Interior.Color = $YELLOW
    }
}

# Write + format one "test case" row: A=name(left cell), B/C/D=field labels,
# E=name repeated (right cell). A gets the left-open box, E the full boxed +
# filled cell, B-D a plain full box.
function Write-Row {
    param($ws, [int]$row, [string]$name, [string]$userName, [string]$passWord, [string]$thirdLabel, [string]$thirdValue)

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $name
    Format-Cell $eCell $false $true

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $name
    Format-Cell $aCell $true $true

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $userName
    Format-Cell $bCell $false $false

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $passWord
    Format-Cell $cCell $false $false

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $thirdLabel
    Format-Cell $dCell $false $false

    $row2 = $row + 1
    $bCell2 = $ws.Cells.Item($row2, 2)
    $bCell2.Value = "Admin"
    Format-Cell $bCell2 $false $false

    $cCell2 = $ws.Cells.Item($row2, 3)
    $cCell2.Value = "admin123"
    Format-Cell $cCell2 $false $false

    $dCell2 = $ws.Cells.Item($row2, 4)
    $dCell2.Value = "zakir"
    Format-Cell $dCell2 $false $false
}

# ---------------------------------------------------------------------------
# Admin sheet
# ---------------------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("Admin")
Write-Row $wsAdmin 1 "Admin_TestCase01" "userName" "passWord" "searchUser" $null
$wsAdmin.Columns.Item(1).AutoFit()
$wsAdmin.Columns.Item(2).AutoFit()
$wsAdmin.Columns.Item(3).AutoFit()
$wsAdmin.Columns.Item(4).AutoFit()
$wsAdmin.Columns.Item(5).AutoFit()
$wsAdmin.Range("E11").Select()

# ---------------------------------------------------------------------------
# PIM sheet
# ---------------------------------------------------------------------------
$wsPim = $wb.Worksheets.Item("PIM")
$wsPim.Range("A2").ClearContents()
Write-Row $wsPim 1 "PIM_TestCase01" "userName" "passWord" "searchUser" $null
$wsPim.Columns.Item(1).AutoFit()
$wsPim.Columns.Item(5).AutoFit()
$wsPim.Range("A1:E2").Select()

# ---------------------------------------------------------------------------
# New TIME sheet
# ---------------------------------------------------------------------------
$wsTime = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTime.Name = "TIME"

Write-Row $wsTime 1 "aa" "userName" "passWord" "empName" $null
Write-Row $wsTime 3 "bb" "userName" "passWord" "empName" $null
Write-Row $wsTime 5 "Time_TestCase01" "userName" "passWord" "empName" $null
Write-Row $wsTime 7 "cc" "userName" "passWord" "empName" $null

# Extra logged columns (E/F) only present for the "Time_TestCase01" block.
$e5 = $wsTime.Cells.Item(5, 5)
$e5.Value = "test"
Format-Cell $e5 $false $false
$f5 = $wsTime.Cells.Item(5, 6)
$f5.Value = "ere"
Format-Cell $f5 $false $false
$g5 = $wsTime.Cells.Item(5, 7)
$g5.Value = "Time_TestCase01"
Format-Cell $g5 $false $true

$e6 = $wsTime.Cells.Item(6, 5)
$e6.Value = "erere"
Format-Cell $e6 $false $false
$f6 = $wsTime.Cells.Item(6, 6)
$f6.Value = "eress"
Format-Cell $f6 $false $false

$wsTime.Columns.Item(1).AutoFit()
$wsTime.Columns.Item(2).AutoFit()
$wsTime.Columns.Item(3).AutoFit()
$wsTime.Columns.Item(4).AutoFit()
$wsTime.Columns.Item(5).AutoFit()
$wsTime.Columns.Item(6).AutoFit()
$wsTime.Columns.Item(7).AutoFit()

$wsTime.Activate()
$wsTime.Range("F11").Select()

Write-Output "done"
